$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "iaest-measure:..." -> "iaest-dimension:..." for D3/E3
$ws.Range("D3").Value = "iaest-dimension:sector-vab-descripcion"
$ws.Range("E3").Value = "iaest-dimension:subseccion-descripcion"

# Row 4: "medida" -> "dim" for D4/E4
$ws.Range("D4").Value = "dim"
$ws.Range("E4").Value = "dim"

# Row 5: "xsd:string" -> "skos:Concept" for D5/E5
$ws.Range("D5").Value = "skos:Concept"
$ws.Range("E5").Value = "skos:Concept"

# New row 6 with mapping file references
$ws.Range("D6").Value = "mapping-sector-vab-descripcion.xlsx"
$ws.Range("E6").Value = "mapping-subseccion-descripcion.xlsx"

# Match the style used by the rest of the data cells (style index 1)
$ws.Range("D5").Copy()
$ws.Range("D6:E6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
